$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; this pushes the existing header row (row 1,
# with its bold/centered/bordered styling) down to row 2, and every data row
# shifts down by one as well.
$ws.Range("A1").EntireRow.Insert()

# Move the old header's look up onto the new row 1 by copying its format
# (still sitting on row 2 after the shift) before we touch anything else.
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A1:L1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 1 with a numeric index header (0-based column index).
for ($i = 0; $i -lt 12; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $i
}

# The old header row (now row 2) reverts to plain formatting.
$ws.Range("A2:L2").ClearFormats()

# The old header row no longer carries thread_size / material_surface
# labels in columns K and L - clear those two cells (keep them present but
# empty, matching column I which was already blank).
$ws.Range("K2").Value = ""
$ws.Range("K2").NumberFormat = "@"
$ws.Range("L2").Value = ""
$ws.Range("L2").NumberFormat = "@"
